# Update COVID-19 country stats (as of 15 June 2020, 12:55) and reorder a
# few countries whose totals now put them in a different rank position.
#
# - Refresh the "datos actualizados" timestamp.
# - Update totals for Estados Unidos, Banglades, Kuwait, Suiza, Rumania,
#   Nepal, Sri Lanka, Eslovenia, Hong Kong, Burkina Faso and Gibraltar.
# - Iran's case count grew enough to overtake Alemania (row 12 becomes
#   Iran with new numbers, row 13 becomes Alemania with its old numbers).
# - The same "rank swap" pattern happens for Senegal/Uzbekistan (75/76),
#   Bosnia y Herzegovina/Venezuela (93/94), Madagascar/Paraguay (119/120)
#   and Groenlandia/Islas Malvinas, Islas Turcas y Caicos/Santa Sede
#   (206/207, 208/209).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 15 de Junio de 2020 a las 12:55'

# Row 4
$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 2162261
$ws.Cells.Item(4, 3).Value = 33
$ws.Cells.Item(4, 4).Value = 870076
$ws.Cells.Item(4, 5).Value = 1174327
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 117858

# Row 12
$ws.Cells.Item(12, 1).Value = 'Iran'
$ws.Cells.Item(12, 2).Value = 189876
$ws.Cells.Item(12, 3).Value = 2449
$ws.Cells.Item(12, 4).Value = 150590
$ws.Cells.Item(12, 5).Value = 30336
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 113
$ws.Cells.Item(12, 8).Value = 8950

# Row 13
$ws.Cells.Item(13, 1).Value = 'Alemania'
$ws.Cells.Item(13, 2).Value = 187671
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 172600
$ws.Cells.Item(13, 5).Value = 6201
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 8870

# Row 21
$ws.Cells.Item(21, 1).Value = 'Banglades'
$ws.Cells.Item(21, 2).Value = 90619
$ws.Cells.Item(21, 3).Value = 3099
$ws.Cells.Item(21, 4).Value = 18731
$ws.Cells.Item(21, 5).Value = 70679
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 38
$ws.Cells.Item(21, 8).Value = 1209

# Row 36
$ws.Cells.Item(36, 1).Value = 'Kuwait'
$ws.Cells.Item(36, 2).Value = 36431
$ws.Cells.Item(36, 3).Value = 511
$ws.Cells.Item(36, 4).Value = 27531
$ws.Cells.Item(36, 5).Value = 8602
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 2
$ws.Cells.Item(36, 8).Value = 298

# Row 39
$ws.Cells.Item(39, 1).Value = 'Suiza'
$ws.Cells.Item(39, 2).Value = 31131
$ws.Cells.Item(39, 3).Value = 14
$ws.Cells.Item(39, 4).Value = 28800
$ws.Cells.Item(39, 5).Value = 393
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 1938

# Row 46
$ws.Cells.Item(46, 1).Value = 'Rumania'
$ws.Cells.Item(46, 2).Value = 22165
$ws.Cells.Item(46, 3).Value = 166
$ws.Cells.Item(46, 4).Value = 15817
$ws.Cells.Item(46, 5).Value = 4921
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 17
$ws.Cells.Item(46, 8).Value = 1427

# Row 74
$ws.Cells.Item(74, 1).Value = 'Nepal'
$ws.Cells.Item(74, 2).Value = 6211
$ws.Cells.Item(74, 3).Value = 451
$ws.Cells.Item(74, 4).Value = 1041
$ws.Cells.Item(74, 5).Value = 5151
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 19

# Row 75
$ws.Cells.Item(75, 1).Value = 'Senegal'
$ws.Cells.Item(75, 2).Value = 5173
$ws.Cells.Item(75, 3).Value = 83
$ws.Cells.Item(75, 4).Value = 3424
$ws.Cells.Item(75, 5).Value = 1685
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 4
$ws.Cells.Item(75, 8).Value = 64

# Row 76
$ws.Cells.Item(76, 1).Value = 'Uzbekistan'
$ws.Cells.Item(76, 2).Value = 5103
$ws.Cells.Item(76, 3).Value = 23
$ws.Cells.Item(76, 4).Value = 3985
$ws.Cells.Item(76, 5).Value = 1099
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 19

# Row 93
$ws.Cells.Item(93, 1).Value = 'Bosnia y Herzegovina'
$ws.Cells.Item(93, 2).Value = 3040
$ws.Cells.Item(93, 3).Value = 147
$ws.Cells.Item(93, 4).Value = 2145
$ws.Cells.Item(93, 5).Value = 731
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 8).Value = 164

# Row 94
$ws.Cells.Item(94, 1).Value = 'Venezuela'
$ws.Cells.Item(94, 2).Value = 2978
$ws.Cells.Item(94, 3).Value = 74
$ws.Cells.Item(94, 4).Value = 835
$ws.Cells.Item(94, 5).Value = 2118
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 1
$ws.Cells.Item(94, 8).Value = 25

# Row 103
$ws.Cells.Item(103, 1).Value = 'Sri Lanka'
$ws.Cells.Item(103, 2).Value = 1896
$ws.Cells.Item(103, 3).Value = 7
$ws.Cells.Item(103, 4).Value = 1342
$ws.Cells.Item(103, 5).Value = 543
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 11

# Row 113
$ws.Cells.Item(113, 1).Value = 'Eslovenia'
$ws.Cells.Item(113, 2).Value = 1496
$ws.Cells.Item(113, 3).Value = 1
$ws.Cells.Item(113, 4).Value = 1359
$ws.Cells.Item(113, 5).Value = 28
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 109

# Row 119
$ws.Cells.Item(119, 1).Value = 'Madagascar'
$ws.Cells.Item(119, 2).Value = 1290
$ws.Cells.Item(119, 3).Value = 18
$ws.Cells.Item(119, 4).Value = 384
$ws.Cells.Item(119, 5).Value = 896
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 10

# Row 120
$ws.Cells.Item(120, 1).Value = 'Paraguay'
$ws.Cells.Item(120, 2).Value = 1289
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(120, 4).Value = 650
$ws.Cells.Item(120, 5).Value = 628
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 11

# Row 123
$ws.Cells.Item(123, 1).Value = 'Hong Kong'
$ws.Cells.Item(123, 2).Value = 1113
$ws.Cells.Item(123, 3).Value = 3
$ws.Cells.Item(123, 4).Value = 1067
$ws.Cells.Item(123, 5).Value = 42
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 4

# Row 129
$ws.Cells.Item(129, 1).Value = 'Burkina Faso'
$ws.Cells.Item(129, 2).Value = 894
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 804
$ws.Cells.Item(129, 5).Value = 37
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 53

# Row 168
$ws.Cells.Item(168, 1).Value = 'Gibraltar'
$ws.Cells.Item(168, 2).Value = 176
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 174
$ws.Cells.Item(168, 5).Value = 2
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 0

# Row 206
$ws.Cells.Item(206, 1).Value = 'Groenlandia'
$ws.Cells.Item(206, 2).Value = 13
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 13
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

# Row 207
$ws.Cells.Item(207, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(207, 2).Value = 13
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 13
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

# Row 208
$ws.Cells.Item(208, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(208, 2).Value = 12
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 1

# Row 209
$ws.Cells.Item(209, 1).Value = 'Santa Sede'
$ws.Cells.Item(209, 2).Value = 12
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 12
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

